# The presentation currently uses the "Integral" (Red Violet) design on its
# slide master / every slide. The author's edit switches the presentation
# back to the default "Office Theme" colour palette (the palette that used
# to live, unused, in the notes-master's theme part).
#
# PowerPoint exposes the live theme of a slide through
# Slide.ThemeColorScheme -> ThemeColorScheme.Colors(index).RGB, one RGBColor
# slot per theme colour (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -
# the same 1..12 ordering used by the OOXML <a:clrScheme> element). Setting
# each slot's .RGB re-colours the active design for the whole deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette: the built-in "Office" colour scheme (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), expressed as COM RGB() long values
# (0xBBGGRR, i.e. red + green*256 + blue*65536).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
